# Add a new column AE with data for date 16-10-2020, following the same
# pattern as the existing AD column (header row + 35 data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous header cell (AD1) onto the new
# header cell (AE1), then set its value to the new date label.
$ws.Range("AD1").Copy($ws.Range("AE1"))
$ws.Range("AE1").Value = "16-10-2020"

# Fill in the data values for the new column, row by row.
$ws.Range("AE2").Value = 55
$ws.Range("AE3").Value = 6357
$ws.Range("AE4").Value = 30
$ws.Range("AE5").Value = 843
$ws.Range("AE6").Value = 972
$ws.Range("AE7").Value = 201
$ws.Range("AE8").Value = 1385
$ws.Range("AE9").Value = 2
$ws.Range("AE10").Value = 5924
$ws.Range("AE11").Value = 525
$ws.Range("AE12").Value = 3606
$ws.Range("AE13").Value = 1623
$ws.Range("AE14").Value = 260
$ws.Range("AE15").Value = 1358
$ws.Range("AE16").Value = 820
$ws.Range("AE17").Value = 10283
$ws.Range("AE18").Value = 1089
$ws.Range("AE19").Value = 65
$ws.Range("AE20").Value = 2710
$ws.Range("AE21").Value = 41196
$ws.Range("AE22").Value = 104
$ws.Range("AE23").Value = 73
$ws.Range("AE24").Value = 0
$ws.Range("AE25").Value = 22
$ws.Range("AE26").Value = 1089
$ws.Range("AE27").Value = 570
$ws.Range("AE28").Value = 3954
$ws.Range("AE29").Value = 1708
$ws.Range("AE30").Value = 59
$ws.Range("AE31").Value = 10472
$ws.Range("AE32").Value = 1256
$ws.Range("AE33").Value = 323
$ws.Range("AE34").Value = 814
$ws.Range("AE35").Value = 6543
$ws.Range("AE36").Value = 5870
